$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-15) to the new homogenization temperature series
$ws.Range("A2").Value = 30
$ws.Range("A3").Value = 25
$ws.Range("A4").Value = 20
$ws.Range("A5").Value = 15
$ws.Range("A6").Value = 10
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 0
$ws.Range("A9").Value = -5
$ws.Range("A10").Value = -10
$ws.Range("A11").Value = -15
$ws.Range("A12").Value = -20
$ws.Range("A13").Value = -25
$ws.Range("A14").Value = -30
$ws.Range("A15").Value = -35

# Scroll back to top-left and move selection to A13, matching final view state
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("A13").Select()
